$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (D1): Schottky diode -> Diode, new part numbers
# (leading apostrophe preserves the original quote-prefixed text style of these cells)
$ws.Range("A12").Value = "'MBRA340T3G"
$ws.Range("B12").Value = "'Diode"
$ws.Range("D12").Value = "'DIOM5226X220N"
$ws.Range("E12").Value = "'MBRA340T3G"

# Row 16 (L1 inductor): new footprint/libref
$ws.Range("D16").Value = "'CDRH8D43NP4R7NC"
$ws.Range("E16").Value = "'CDRH8D43NP-4R7NC"

# Row 17 (L2 inductor): new footprint/libref (both equal to the new part number)
$ws.Range("D17").Value = "'74438357010"
$ws.Range("E17").Value = "'74438357010"

# Uniform column widths A:F as set by the author (~18.8 chars, matches stored width 18.8164)
$ws.Range("A:F").ColumnWidth = 17.983072916666668
